# "Update report and status"
#
# - The "Performance Metrics Appendix" row (row 7) is removed entirely;
#   everything below it shifts up one row.
# - A number of sections progress from Red ("no progress") to either
#   Yellow ("some progress") or Green ("completed").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Performance Metrics Appendix" row - the whole row shifts up.
$ws.Rows.Item(7).Delete()

# Colour constants (VBA-style RGB packed as R + G*256 + B*65536),
# matching the workbook's existing Red / Yellow / Green status fills.
$Red    = 255        # FFFF0000
$Yellow = 65535       # FFFFFF00
$Green  = 5287936     # FF00B050

# Rows below are the *post-delete* row numbers (sections shifted up by one).
$ws.Range("A2").Interior.Color  = $Green    # Challenge Definition: Red -> Green
$ws.Range("A3").Interior.Color  = $Yellow   # Customer Value Proposition: Red -> Yellow
$ws.Range("A4").Interior.Color  = $Yellow   # Changes from Proposal: Red -> Yellow
$ws.Range("A5").Interior.Color  = $Yellow   # Key Technical Elements: Red -> Yellow
$ws.Range("A6").Interior.Color  = $Yellow   # Performance Metrics Summay: Red -> Yellow
$ws.Range("A7").Interior.Color  = $Yellow   # Failure Analysis: Red -> Yellow
$ws.Range("A8").Interior.Color  = $Yellow   # Timeline: Red -> Yellow
$ws.Range("A15").Interior.Color = $Green    # Acknowledgements: Red -> Green
$ws.Range("A17").Interior.Color = $Yellow   # Testing instructions: Red -> Yellow
$ws.Range("A18").Interior.Color = $Yellow   # Assembly instructions: Red -> Yellow
$ws.Range("A19").Interior.Color = $Yellow   # Technical Documentation: Red -> Yellow
